# Error Calculations and Plots
# Remove the "RM 232" and "SC 92" data rows entirely (rows 26 and 28 in the
# original sheet), which shifts all subsequent rows up. Then update a
# handful of F (and a couple of D) column values that changed as part of
# this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (originally row 26). This shifts "SC 92"
# (originally row 28) up to row 27.
$ws.Rows.Item(26).Delete()

# Delete the "SC 92" row, now at row 27 after the first deletion.
$ws.Rows.Item(27).Delete()

# Apply the remaining value changes (positions are post-deletion row numbers).
$ws.Range("F19").Value() = 17.81
$ws.Range("F21").ClearContents()
$ws.Range("F23").Value() = 16.48
$ws.Range("D26").ClearContents()
$ws.Range("D27").Value() = -14.6
$ws.Range("F27").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("F33").Value() = 17.53
